$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Row 3 ("Added ability to make friend,battle,breed requests" row),
#        3rd cell: " requests" -> " requests " (the run gains a trailing
#        space) and the old redundant trailing run (a lone space that used
#        to sit after the _GoBack bookmark) goes away, since the bookmark
#        is being relocated to the end of the new row below (see step 2).
$reqCell = $t.Cell(3, 3)
$reqRange = $reqCell.Range
$null = $reqRange.Find.Execute(" requests", $false, $false, $false, $false, $false, $true, 1, $false, " requests ", 2)

# Remove the trailing "s " (last letter of "requests" + the old lone-space
# run) and retype it so the cell ends up with a single clean run again and
# the stale _GoBack bookmark that used to live here is cleared out (it gets
# re-created at its new home in step 2).
$cellEnd = $t.Cell(3, 3).Range.End
$tailSpan = $d.Range($cellEnd - 3, $cellEnd - 1)
$tailSpan.Delete()
$cellEnd2 = $t.Cell(3, 3).Range.End
$d.Range($cellEnd2 - 1, $cellEnd2 - 1).InsertBefore("s ")

# --- 2. Fill in the next (previously blank) row with the new timesheet
#        entry, and park the _GoBack bookmark at the end of its last cell
#        (where Word leaves it after the most recent edit/typing).
$t.Cell(4, 1).Range.Text = "27/11/2012"
$t.Cell(4, 2).Range.Text = "6pm-8pm"
$t.Cell(4, 3).Range.Text = "Server to server documentation"

$newBmRange = $t.Cell(4, 3).Range
$newBmRange.Collapse(0)
$null = $d.Bookmarks.Add("_GoBack", $newBmRange)
